# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# figures for the coin rows, matching the upstream GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new display text.
$updates = @(
    @{ Cell = 'D2'; Value = '63.363.12' }
    @{ Cell = 'E2'; Value = '  +0.65%  ' }
    @{ Cell = 'D3'; Value = '2.677.65' }
    @{ Cell = 'E3'; Value = '  +3.94%  ' }
    @{ Cell = 'E4'; Value = '  +0.07%  ' }
    @{ Cell = 'D5'; Value = '610.70' }
    @{ Cell = 'E5'; Value = '  +4.43%  ' }
    @{ Cell = 'D6'; Value = '143.67' }
    @{ Cell = 'E6'; Value = '  -0.47%  ' }
    @{ Cell = 'E7'; Value = '  +0.08%  ' }
    @{ Cell = 'D8'; Value = '0.587' }
    @{ Cell = 'E8'; Value = '  -0.54%  ' }
    @{ Cell = 'D9'; Value = '2.678.95' }
    @{ Cell = 'E9'; Value = '  +4.00%  ' }
    @{ Cell = 'E10'; Value = '  +0.74%  ' }
    @{ Cell = 'D11'; Value = '5.62' }
    @{ Cell = 'E11'; Value = '  +0.63%  ' }
    @{ Cell = 'E12'; Value = '  +0.60%  ' }
    @{ Cell = 'D13'; Value = '0.361' }
    @{ Cell = 'E13'; Value = '  +3.23%  ' }
    @{ Cell = 'D14'; Value = '27.35' }
    @{ Cell = 'E14'; Value = '  +1.07%  ' }
    @{ Cell = 'D15'; Value = '3.161.17' }
    @{ Cell = 'E15'; Value = '  +4.00%  ' }
    @{ Cell = 'D16'; Value = '63.309.38' }
    @{ Cell = 'D17'; Value = '0.0000145' }
    @{ Cell = 'E17'; Value = '  +0.29%  ' }
    @{ Cell = 'D18'; Value = '2.690.06' }
    @{ Cell = 'E18'; Value = '  +4.63%  ' }
    @{ Cell = 'D19'; Value = '11.46' }
    @{ Cell = 'E19'; Value = '  +3.48%  ' }
    @{ Cell = 'D20'; Value = '342.61' }
    @{ Cell = 'E20'; Value = '  +0.33%  ' }
    @{ Cell = 'D21'; Value = '4.41' }
    @{ Cell = 'E21'; Value = '  +1.80%  ' }
    @{ Cell = 'D22'; Value = '6.87' }
    @{ Cell = 'E22'; Value = '  +3.67%  ' }
    @{ Cell = 'E23'; Value = '  +0.15%  ' }
    @{ Cell = 'D24'; Value = '67.45' }
    @{ Cell = 'E24'; Value = '  -0.48%  ' }
    @{ Cell = 'D25'; Value = '1.65' }
    @{ Cell = 'E25'; Value = '  +3.03%  ' }
    @{ Cell = 'D26'; Value = '1.54' }
    @{ Cell = 'E26'; Value = '  -3.36%  ' }
    @{ Cell = 'D27'; Value = '8.64' }
    @{ Cell = 'E27'; Value = '  +4.84%  ' }
    @{ Cell = 'D28'; Value = '0.164' }
    @{ Cell = 'E28'; Value = '  -0.44%  ' }
    @{ Cell = 'D29'; Value = '543.26' }
    @{ Cell = 'E29'; Value = '  +18.12%  ' }
    @{ Cell = 'E30'; Value = '  +0.12%  ' }
    @{ Cell = 'D31'; Value = '7.88' }
    @{ Cell = 'E31'; Value = '  -1.54%  ' }
    @{ Cell = 'E32'; Value = '  +6.87%  ' }
    @{ Cell = 'E33'; Value = '  +7.96%  ' }
    @{ Cell = 'D34'; Value = '0.0₃0809' }
    @{ Cell = 'E34'; Value = '  +1.43%  ' }
    @{ Cell = 'D35'; Value = '172.66' }
    @{ Cell = 'E35'; Value = '  -2.34%  ' }
    @{ Cell = 'D36'; Value = '5.15' }
    @{ Cell = 'E36'; Value = '  +14.05%  ' }
    @{ Cell = 'D37'; Value = '0.405' }
    @{ Cell = 'E37'; Value = '  +1.77%  ' }
    @{ Cell = 'E38'; Value = '  -0.03%  ' }
    @{ Cell = 'D39'; Value = '19.26' }
    @{ Cell = 'E39'; Value = '  +2.27%  ' }
    @{ Cell = 'D40'; Value = '1.85' }
    @{ Cell = 'E40'; Value = '  +9.28%  ' }
    @{ Cell = 'D41'; Value = '176.59' }
    @{ Cell = 'E41'; Value = '  +11.60%  ' }
    @{ Cell = 'D42'; Value = '0.998' }
    @{ Cell = 'E42'; Value = '  -0.10%  ' }
    @{ Cell = 'D43'; Value = '3.75' }
    @{ Cell = 'E43'; Value = '  +1.82%  ' }
    @{ Cell = 'D44'; Value = '22.19' }
    @{ Cell = 'E44'; Value = '  +4.56%  ' }
    @{ Cell = 'D45'; Value = '0.0569' }
    @{ Cell = 'E45'; Value = '  +6.14%  ' }
    @{ Cell = 'D46'; Value = '0.636' }
    @{ Cell = 'E46'; Value = '  +0.17%  ' }
    @{ Cell = 'D47'; Value = '0.0966' }
    @{ Cell = 'E47'; Value = '  +0.71%  ' }
    @{ Cell = 'D48'; Value = '0.0240' }
    @{ Cell = 'E48'; Value = '  +2.18%  ' }
    @{ Cell = 'D49'; Value = '18.98' }
    @{ Cell = 'E49'; Value = '  +5.22%  ' }
    @{ Cell = 'D50'; Value = '1.75' }
    @{ Cell = 'E50'; Value = '  +4.25%  ' }
    @{ Cell = 'D51'; Value = '11.29' }
    @{ Cell = 'E51'; Value = '  -0.99%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Columns D/E store these as text (e.g. '610.70', '0.0240', '  +4.43%  ').
    # Excel auto-converts plain numeric-looking text to a Number on assignment,
    # which would drop formatting like trailing zeros - force Text first, then
    # restore the default 'Normal' style so no stray numeric format lingers.
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
    $rng.Style = 'Normal'
}

